# Applies the cryptocurrency price/volume refresh described by the commit
# "Updated cryptos list on Sat Oct 19 21:32:40 UTC 2024 with GitHub Actions".
# Column D (Price) and Column E (Volume(1h)) values are plain text in the
# sheet (e.g. "597.50", "  -0.12%  "), so each target cell is forced to the
# Text number format before the new value is written — this stops Excel's
# COM layer from re-interpreting strings like "597.50" or "0.142" as
# numbers (which would silently drop the trailing zero / introduce binary
# floating point noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.333.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.649.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.69'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.142'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.09%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.352'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.11'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.131.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.235.79'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.654.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.40'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.33'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.82%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '559.50'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.07'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.86'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.73'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.35'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₆0335'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.62'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '159.02'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.20'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.615'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.568'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.03%  '
